$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1049.8823
$ws.Range("I11").Value = 1049.8823
$ws.Range("K11").Value = 1049.8823
$ws.Range("M11").Value = -909.8823

$ws.Range("H15").Value = 1275.7906
$ws.Range("I15").Value = 1275.7906
$ws.Range("K15").Value = 3827.3718
$ws.Range("M15").Value = -3658.3718

$ws.Range("H40").Value = 2117.92
$ws.Range("I40").Value = 2791.4167
$ws.Range("J40").Value = 1496.2307
$ws.Range("K40").Value = 2791.4167
$ws.Range("L40").Value = 1496.2307
$ws.Range("M40").Value = -2616.4167
$ws.Range("N40").Value = -1846.2307

$ws.Range("H53").Value = 259.5
$ws.Range("I53").Value = 138.36363
$ws.Range("J53").Value = 407.55554
$ws.Range("K53").Value = 138.36363
$ws.Range("L53").Value = 407.55554
$ws.Range("M53").Value = 498.63637
$ws.Range("N53").Value = -1681.55554

$ws.Range("H69").Value = 4253.25
$ws.Range("J69").Value = 4253.25
$ws.Range("L69").Value = 12759.75
$ws.Range("N69").Value = -14507.75

$ws.Range("H72").Value = 4253.25
$ws.Range("J72").Value = 4253.25
$ws.Range("L72").Value = 38279.25
$ws.Range("N72").Value = -47015.25

$ws.Range("H74").Value = 4653.231
$ws.Range("I74").Value = 4013.2856
$ws.Range("K74").Value = 4013.2856
$ws.Range("M74").Value = -3077.2856

$ws.Range("H77").Value = 4653.231
$ws.Range("I77").Value = 4013.2856
$ws.Range("K77").Value = 20066.428
$ws.Range("M77").Value = -15386.428

$ws.Range("H113").Value = 2679.5
$ws.Range("J113").Value = 2906
$ws.Range("L113").Value = 2906
$ws.Range("N113").Value = -9414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6282.6177
$ws.Range("I61").Value = 4925.7407
$ws.Range("K61").Value = 4925.7407
$ws.Range("M61").Value = -4713.7407

$ws.Range("H74").Value = 2576.3794
$ws.Range("I74").Value = 2227.158
$ws.Range("J74").Value = 3239.9
$ws.Range("K74").Value = 2227.158
$ws.Range("L74").Value = 3239.9
$ws.Range("M74").Value = -1353.158
$ws.Range("N74").Value = -4987.9

$ws.Range("H77").Value = 2576.3794
$ws.Range("I77").Value = 2227.158
$ws.Range("J77").Value = 3239.9
$ws.Range("K77").Value = 11135.79
$ws.Range("L77").Value = 16199.5
$ws.Range("M77").Value = -6767.789999999999
$ws.Range("N77").Value = -24935.5

$ws.Range("H102").Value = 3255.3333
$ws.Range("I102").Value = 2749.6667
$ws.Range("J102").Value = 4266.6665
$ws.Range("K102").Value = 2749.6667
$ws.Range("L102").Value = 4266.6665
$ws.Range("M102").Value = -1127.6667
$ws.Range("N102").Value = -7510.6665

$ws.Range("H136").Value = 6282.6177
$ws.Range("I136").Value = 4925.7407
$ws.Range("K136").Value = 14777.2221
$ws.Range("M136").Value = -12227.2221

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H94").Value = 941
$ws.Range("I94").Value = 1375
$ws.Range("J94").Value = 651.6667
$ws.Range("K94").Value = 1375
$ws.Range("L94").Value = 651.6667
$ws.Range("M94").Value = -924
$ws.Range("N94").Value = -1553.6667

$ws.Range("H99").Value = 1013.2143
$ws.Range("I99").Value = 932.1667
$ws.Range("J99").Value = 1499.5
$ws.Range("K99").Value = 932.1667
$ws.Range("L99").Value = 1499.5
$ws.Range("M99").Value = 565.8333
$ws.Range("N99").Value = -4495.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13007.261
$ws.Range("I31").Value = 12686.294
$ws.Range("J31").Value = 13916.667
$ws.Range("K31").Value = 12686.294
$ws.Range("L31").Value = 13916.667
$ws.Range("M31").Value = -12391.294
$ws.Range("N31").Value = -14506.667

$ws.Range("H34").Value = 13007.261
$ws.Range("I34").Value = 12686.294
$ws.Range("J34").Value = 13916.667
$ws.Range("K34").Value = 12686.294
$ws.Range("L34").Value = 13916.667
$ws.Range("M34").Value = -12484.294
$ws.Range("N34").Value = -14320.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 28000
$ws.Range("J74").Value = 28000
$ws.Range("L74").Value = 28000
$ws.Range("N74").Value = -29872

$ws.Range("H77").Value = 28000
$ws.Range("J77").Value = 28000
$ws.Range("L77").Value = 84000
$ws.Range("N77").Value = -93360

$ws.Range("H80").Value = 8320
$ws.Range("I80").Value = 3600
$ws.Range("J80").Value = 9500
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 9500
$ws.Range("M80").Value = -2602
$ws.Range("N80").Value = -11496

$ws.Range("H83").Value = 8320
$ws.Range("I83").Value = 3600
$ws.Range("J83").Value = 9500
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 47500
$ws.Range("M83").Value = -13008
$ws.Range("N83").Value = -57484

$ws.Range("H113").Value = 2264.2856
$ws.Range("I113").Value = 2172.7273
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 2172.7273
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = -2.727300000000014
$ws.Range("N113").Value = -6940

$ws.Range("H134").Value = 23116.428
$ws.Range("J134").Value = 23116.428
$ws.Range("L134").Value = 69349.284
$ws.Range("N134").Value = -74419.284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000.5
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H82").Value = 1606.6
$ws.Range("J82").Value = 1140.6
$ws.Range("L82").Value = 1140.6
$ws.Range("N82").Value = -1862.6

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H85").Value = 1606.6
$ws.Range("J85").Value = 1140.6
$ws.Range("L85").Value = 1140.6
$ws.Range("N85").Value = -3636.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 64808.332
$ws.Range("J133").Value = 64808.332
$ws.Range("L133").Value = 64808.332
$ws.Range("N133").Value = -74928.33199999999

$ws.Range("H136").Value = 3838.3696
$ws.Range("I136").Value = 3060.8518
$ws.Range("J136").Value = 4943.263
$ws.Range("K136").Value = 9182.555399999999
$ws.Range("L136").Value = 14829.789
$ws.Range("M136").Value = -6632.555399999999
$ws.Range("N136").Value = -19929.789
